$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 298, shifting existing data (rows 298-320) down to 301-323.
$ws.Rows.Item(298).Resize(3).Insert()

# New row 298: Chirimoya - Especial, 2023-10-13
$ws.Cells.Item(298, 1).Value = 8
$ws.Cells.Item(298, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(298, 3).Value = "Coquimbo"
$ws.Cells.Item(298, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(298, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(298, 5).Value = 4
$ws.Cells.Item(298, 6).Value = "Fruta"
$ws.Cells.Item(298, 7).Value = 100107
$ws.Cells.Item(298, 8).Value = "Otros"
$ws.Cells.Item(298, 9).Value = 100107002
$ws.Cells.Item(298, 10).Value = "Chirimoya"
$ws.Cells.Item(298, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(298, 12).Value = "Especial"
$ws.Cells.Item(298, 13).Value = 160
$ws.Cells.Item(298, 14).Value = 20000
$ws.Cells.Item(298, 15).Value = 21000
$ws.Cells.Item(298, 16).Value = 20500
$ws.Cells.Item(298, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(298, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(298, 19).Value = 2050
$ws.Cells.Item(298, 20).Value = 10

# New row 299: Chirimoya - Primera, 2023-10-13
$ws.Cells.Item(299, 1).Value = 8
$ws.Cells.Item(299, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(299, 3).Value = "Coquimbo"
$ws.Cells.Item(299, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(299, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(299, 5).Value = 4
$ws.Cells.Item(299, 6).Value = "Fruta"
$ws.Cells.Item(299, 7).Value = 100107
$ws.Cells.Item(299, 8).Value = "Otros"
$ws.Cells.Item(299, 9).Value = 100107002
$ws.Cells.Item(299, 10).Value = "Chirimoya"
$ws.Cells.Item(299, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(299, 12).Value = "Primera"
$ws.Cells.Item(299, 13).Value = 200
$ws.Cells.Item(299, 14).Value = 17000
$ws.Cells.Item(299, 15).Value = 18000
$ws.Cells.Item(299, 16).Value = 17500
$ws.Cells.Item(299, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(299, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(299, 19).Value = 1750
$ws.Cells.Item(299, 20).Value = 10

# New row 300: Chirimoya - Segunda, 2023-10-13
$ws.Cells.Item(300, 1).Value = 8
$ws.Cells.Item(300, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(300, 3).Value = "Coquimbo"
$ws.Cells.Item(300, 4).Value = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(300, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(300, 5).Value = 4
$ws.Cells.Item(300, 6).Value = "Fruta"
$ws.Cells.Item(300, 7).Value = 100107
$ws.Cells.Item(300, 8).Value = "Otros"
$ws.Cells.Item(300, 9).Value = 100107002
$ws.Cells.Item(300, 10).Value = "Chirimoya"
$ws.Cells.Item(300, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(300, 12).Value = "Segunda"
$ws.Cells.Item(300, 13).Value = 200
$ws.Cells.Item(300, 14).Value = 13000
$ws.Cells.Item(300, 15).Value = 14000
$ws.Cells.Item(300, 16).Value = 13500
$ws.Cells.Item(300, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(300, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(300, 19).Value = 1350
$ws.Cells.Item(300, 20).Value = 10
